# Sprint Backlog.xlsx update
# - Move the existing "Sugerir 2 Features (...)" suggestions from column A
#   (rows 2-6) into column D (rows 2-6).
# - Put the two newly agreed user-story tasks into A2/A3.
# - Clear A4:A6 of text but mark them with the underlined style that is used
#   to show "not started yet" placeholder rows.
# - Add a styled (underlined), empty placeholder cell at B8.
# - Widen column A a bit and move the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "Sugerir 2 Features (...)" values before they are
# overwritten, then re-home them in column D (same row).
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $ws.Cells.Item($r, 4).Value = $cell.Value()
}

# New user-story tasks for the first two rows of column A.
$ws.Range("A2").Value = "Criar User Story Primeira Feature"
$ws.Range("A3").Value = "Criar User Story Segunda Feature"

# Clear the remaining column-A cells and give them the underline style used
# for still-empty backlog rows.
$ws.Range("A4").Value = ""
$ws.Range("A5").Value = ""
$ws.Range("A6").Value = ""
$ws.Range("A4:A6").Font.Underline = $true

# New empty styled placeholder row.
$ws.Range("B8").Font.Underline = $true

# Column A grew a bit wider to fit the new text.
$ws.Columns.Item(1).ColumnWidth = 32.109375

# Leave the user's selection on the first newly-edited cell.
$ws.Range("A3").Select()
